$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PEIIR-EFPEIIR")

# Issues #280 / #99: add new power plant types to the Electricity Source
# subscript used by the Electricity Fuel PEIIR table.
[void]$ws.Activate()

$ws.Range("S1").Value = "hard coal w CCS"
$ws.Range("T1").Value = "natural gas combined cycle w CCS"
$ws.Range("U1").Value = "biomass w CCS"
$ws.Range("V1").Value = "lignite w CCS"
$ws.Range("W1").Value = "small modular reactor"
$ws.Range("X1").Value = "hydrogen"

# Match the header formatting already used across the rest of row 1
# (right aligned, wrapped header text).
$ws.Range("S1:X1").HorizontalAlignment = -4152
$ws.Range("S1:X1").WrapText = $true

# New columns get the same width as the other data columns.
$ws.Range("S1:X13").ColumnWidth = 13.42578125

# Fill the new columns with the same default "0" improvement rate used
# throughout the rest of the table.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 19).Value = 0
    $ws.Cells.Item($r, 20).Value = 0
    $ws.Cells.Item($r, 21).Value = 0
    $ws.Cells.Item($r, 22).Value = 0
    $ws.Cells.Item($r, 23).Value = 0
    $ws.Cells.Item($r, 24).Value = 0
}

# A trailing styled (but still empty) column, matching the extra header
# formatting left just past the newly added data.
$ws.Range("Y1").HorizontalAlignment = -4152
$ws.Range("Y1").WrapText = $true

[void]$ws.Range("X2").Select()

# Return to the workbook's original active sheet.
[void]$wb.Worksheets.Item("About").Activate()
